$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing address value for the first staff row (H2),
# matching the existing "Thủ Đức, TP.HCM" value used elsewhere (H4).
$ws.Range("H2").Value = "Thủ Đức, TP.HCM"

# Reflect final selection state left after the edit.
$ws.Range("H7").Select()
